$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C (rows 2-235) held a fractional depreciation/inflation rate
# (e.g. 0.2637...) that should instead be expressed as a percentage
# value (e.g. 26.37...), so every cell in C2:C235 is rescaled by 100.
$rng = $ws.Range("C2:C235")
foreach ($cell in $rng) {
    $cell.Value = $cell.Value2 * 100
}

# Move the active selection to C4, matching the saved cursor position.
$ws.Range("C4").Select()
